$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pagseguro")

# Set column A (enable) to 1 for rows 2-8
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = 1
}

# Update the selection on the sheet to C19
$ws.Range("C19").Select()
